$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.523.35'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '3.326.40'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '595.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  +2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.423'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '3.903.62'
$ws.Range('E12').Value = '  +2.39%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '69.540.40'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '3.327.97'
$ws.Range('E17').Value = '  +3.11%  '
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '422.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.520'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.193'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '164.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '27.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.804'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('D41').Value = '2.718.05'
$ws.Range('E41').Value = '  +4.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.02%  '
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '345.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('E49').Value = '  +2.82%  '
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('E51').Value = '  -0.52%  '
